$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": narrower columns + simplified header text ---
$resumen = $wb.Worksheets.Item("Resumen")

$resumen.Range("A1").Value = "Administración de costos del proyecto"
$resumen.Range("B1").Value = "Costos de proyecto"

# Column widths (COM ColumnWidth is in characters; Excel persists the OOXML
# <col width> with a +5/MDW padding and snaps to pixel granularity, so we
# pass the character width that lands closest to the target stored width).
$resumen.Columns.Item(1).ColumnWidth = 49.285714285714285
$resumen.Columns.Item(2).ColumnWidth = 30.428571428571427

# --- Sheet "Etapas_proyecto": drop the trailing "Control de cambios" block ---
$etapas = $wb.Worksheets.Item("Etapas_proyecto")
$etapas.Rows("32:35").Delete()
